$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Table layout before the edit: row 1 = header, row 2 = the single data
# lot ("1."), row 3 = the "Razem:" totals row. We insert two more data
# lots ("2." and "3.") directly above the totals row, then refresh the
# totals. Row references in this host resolve by position, so after each
# insertion we re-fetch the row we want to insert in front of by its
# (now shifted) index rather than reusing the old row object.

# --- new lot "2." (inserted right after the existing "1." row) ---
$beforeRow = $t.Rows.Item(3)
$row2 = $t.Rows.Add($beforeRow)
$row2.Cells.Item(1).Range.Text  = "2."
$row2.Cells.Item(2).Range.Text  = "7"
$row2.Cells.Item(3).Range.Text  = "Opoczno"
$row2.Cells.Item(4).Range.Text  = "7891234567"
$row2.Cells.Item(5).Range.Text  = "DB"
$row2.Cells.Item(6).Range.Text  = "C"
$row2.Cells.Item(7).Range.Text  = "2,50"
$row2.Cells.Item(8).Range.Text  = "35,00"
$row2.Cells.Item(9).Range.Text  = "4,00"
$row2.Cells.Item(10).Range.Text = "3 000,00"
$row2.Cells.Item(11).Range.Text = "12 000,00"

# --- new lot "3." (inserted right after the row just added) ---
$beforeRow = $t.Rows.Item(4)
$row3 = $t.Rows.Add($beforeRow)
$row3.Cells.Item(1).Range.Text  = "3."
$row3.Cells.Item(2).Range.Text  = "8"
$row3.Cells.Item(3).Range.Text  = "Poddębice"
$row3.Cells.Item(4).Range.Text  = "8912345678"
$row3.Cells.Item(5).Range.Text  = "DB"
$row3.Cells.Item(6).Range.Text  = "C"
$row3.Cells.Item(7).Range.Text  = "2,50"
$row3.Cells.Item(8).Range.Text  = "35,00"
$row3.Cells.Item(9).Range.Text  = "4,00"
$row3.Cells.Item(10).Range.Text = "1 500,00"
$row3.Cells.Item(11).Range.Text = "6 000,00"

# --- refresh the "Razem:" totals (now the last row). Note the first
# logical cell of this row spans 8 grid columns (gridSpan), so the row
# only has 4 logical cells: "Razem:", Masa total, "Razem:", Wartosc total.
$razemRow = $t.Rows.Item($t.Rows.Count)
$razemRow.Cells.Item(2).Range.Text = "12,00"
$razemRow.Cells.Item(4).Range.Text = "30 000,00"
